# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The sheet tracks per-outing pitching data. Column G (header "K") holds a
# stat that was recomputed from the source data and rewritten here; the
# other columns (C..F, H..J) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values for rows 2..62, in order.
$kValues = @(
    2, 0, 0, 1, 1, 0, 1, 0, 4, 1,
    2, 3, 1, 1, 2, 4, 2, 0, 1, 1,
    0, 2, 2, 3, 3, 0, 0, 1, 1, 0,
    0, 2, 1, 1, 2, 0, 0, 1, 3, 1,
    2, 4, 1, 2, 1, 1, 1, 2, 1, 0,
    0, 0, 0, 2, 0, 1, 2, 1, 1, 2,
    1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
